$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance marks (value 5) for the newly-covered rows
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5

$ws.Range("C18").Value = 5

# Move the active cell/selection on the frozen bottom-right pane to C21
$ws.Range("C21").Select()
